$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new quarterly data point (01-07-2021 / 192.0109) as row 40.
# The date-like text "01-07-2021" must land in the shared-string table as a
# literal string, not get auto-converted to a serial date number. Using a
# text formula that evaluates to the literal string, then collapsing it to
# a static value via copy / paste-special-values, writes a plain `t="s"`
# string cell without perturbing any cell's NumberFormat/style.
$ws.Range("A40").Formula = "=""01-07-2021"""
$ws.Range("A40").Copy()
$ws.Range("A40").PasteSpecial(-4163)

$ws.Range("B40").Value = 192.0109
